# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect a fresh report run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for the first data row
$overview.Range("G2").Value = "2016-08-31 08:10:00"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
$zhcn.Range("H2").Value = "2016-08-31 08:09:49"
$zhcn.Range("K2").Value = "2016-08-31 08:10:46"

# de-de sheet: "Correspond Handoff Datetime" (H)
$dede.Range("H2").Value = "2016-08-31 08:11:06"
